$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 118 (shifts existing rows 118:146 down to 119:147,
# inheriting the formatting of the row above as Excel normally does).
$ws.Rows.Item(118).Insert()

# Populate the new row 118 with the new weekly price-report record (same static
# columns as the rest of the "Jengibre" block, only the date/volume/price columns differ).
$ws.Range("A118").Value = 10
$ws.Range("B118").Value = "Vega Modelo de Temuco"
$ws.Range("C118").Value = "La Araucanía"
$ws.Range("D118").Value = 44642
$ws.Range("E118").Value = 9
$ws.Range("F118").Value = 100114007
$ws.Range("G118").Value = "Jengibre"
$ws.Range("H118").Value = "Sin especificar"
$ws.Range("I118").Value = "Primera"
$ws.Range("J118").Value = 15
$ws.Range("K118").Value = 25000
$ws.Range("L118").Value = 25000
$ws.Range("M118").Value = 25000
$ws.Range("N118").Value = "$/caja 13 kilos"
$ws.Range("O118").Value = "Perú"
$ws.Range("P118").Value = 1923
$ws.Range("Q118").Value = 13
$ws.Range("R118").Value = "Hortaliza"
